# Edit script: insert two new records (rows) into the "Femacal de La Calera - Zapallo"
# hortaliza sheet, shifting the existing rows 237..316 down to 239..318, and
# populating the two newly inserted rows (237 and 238) with fresh data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 237; Excel shifts rows 237-316 down to 239-318
# and copies formatting (e.g. the date number format on column D) from the
# row above, just like a normal Excel "Insert Copied Cells"/"Insert Rows".
$ws.Rows.Item(237).Resize(2).Insert()

# --- Row 237: Camote, "1a (guarda)" ---
$ws.Cells.Item(237, 1).Value = 3
$ws.Cells.Item(237, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(237, 3).Value = "Coquimbo"
$ws.Cells.Item(237, 4).Value2 = 44468
$ws.Cells.Item(237, 5).Value = 5
$ws.Cells.Item(237, 6).Value = 100112045
$ws.Cells.Item(237, 7).Value = "Zapallo"
$ws.Cells.Item(237, 8).Value = "Camote"
$ws.Cells.Item(237, 9).Value = "1a (guarda)"
$ws.Cells.Item(237, 10).Value = 110
$ws.Cells.Item(237, 11).Value = 800
$ws.Cells.Item(237, 12).Value = 800
$ws.Cells.Item(237, 13).Value = 800
$ws.Cells.Item(237, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(237, 15).Value = "Provincia de Talca"
$ws.Cells.Item(237, 16).Value = 800
$ws.Cells.Item(237, 17).Value = 1
$ws.Cells.Item(237, 18).Value = "Hortaliza"

# --- Row 238: Camote, "2a (guarda)" ---
$ws.Cells.Item(238, 1).Value = 3
$ws.Cells.Item(238, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(238, 3).Value = "Coquimbo"
$ws.Cells.Item(238, 4).Value2 = 44468
$ws.Cells.Item(238, 5).Value = 5
$ws.Cells.Item(238, 6).Value = 100112045
$ws.Cells.Item(238, 7).Value = "Zapallo"
$ws.Cells.Item(238, 8).Value = "Camote"
$ws.Cells.Item(238, 9).Value = "2a (guarda)"
$ws.Cells.Item(238, 10).Value = 120
$ws.Cells.Item(238, 11).Value = 600
$ws.Cells.Item(238, 12).Value = 600
$ws.Cells.Item(238, 13).Value = 600
$ws.Cells.Item(238, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(238, 15).Value = "Provincia de Talca"
$ws.Cells.Item(238, 16).Value = 600
$ws.Cells.Item(238, 17).Value = 1
$ws.Cells.Item(238, 18).Value = "Hortaliza"
